$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add emails for the users that didn't have one yet (E6:G6), with hyperlinks,
# copying the existing formatting from the analogous filled-in cells.
$ws.Range("E6").Value = "ahenao@edeq.com"
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:ahenao@edeq.com")
$ws.Range("D6").Copy()
$ws.Range("E6").PasteSpecial(-4122)

$ws.Range("F6").Value = "slopez@gmail.com"
$ws.Hyperlinks.Add($ws.Range("F6"), "mailto:slopez@gmail.com")
$ws.Range("D6").Copy()
$ws.Range("F6").PasteSpecial(-4122)

$ws.Range("G6").Value = "storres@hotmail.com"
$ws.Hyperlinks.Add($ws.Range("G6"), "mailto:storres@hotmail.com")
$ws.Range("B6").Copy()
$ws.Range("G6").PasteSpecial(-4122)

# Fix typo: sgracia -> sgarcia
$ws.Range("H5").Value = "sgarcia"

# Email for the just-renamed user
$ws.Range("H6").Value = "sgarcia@edeq.com"
$ws.Hyperlinks.Add($ws.Range("H6"), "mailto:sgarcia@edeq.com")
$ws.Range("B6").Copy()
$ws.Range("H6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Add two new users (columns I and J) with names
$ws.Range("I3").Value = "Juan David restrepo"
$ws.Range("J3").Value = "Lina maria duran"

# Move the selection as it ended up after the edits
[void]$ws.Range("B11").Select()
